$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13: Recipient name text update
$ws.Range("C13").Value = "Satheesh Nair"

# Row 16: replace the volatile =TODAY() formula with a fixed literal value,
# and switch the cell's number format from a date format to a plain integer.
$ws.Range("C16").Value = 12312020
$ws.Range("C16").NumberFormat = "0"

# Update the active cell / selection on the sheet.
$ws.Range("E17").Select()

# Bump the vertical print resolution recorded on the page setup.
$ws.PageSetup.PrintQuality = 300
